$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Current Status" column (D) for several tasks to reflect
# project progress. We copy/paste-special (formats only) from an existing
# cell that already carries the desired "Done" / "In Progress" look so the
# workbook's existing style records are reused instead of new ones being
# created, then set the text.

# Row moving to "In Progress": Demo Video (copy format from D7 first,
# before D7 itself is updated below)
$ws.Range("D7").Copy()
$ws.Range("D14").PasteSpecial(-4122)
$ws.Range("D14").Value = "In Progress"

# Rows moving to "Done": Interactive Demo, Final Product 1st Release,
# User Manual, Project Report, Test Results, Project Tracking Form,
# Project Poster, Project Webpage Update
$doneRows = @(6, 7, 8, 9, 10, 11, 12, 13)
foreach ($r in $doneRows) {
    $ws.Range("D3").Copy()
    $ws.Range("D$r").PasteSpecial(-4122)
    $ws.Range("D$r").Value = "Done"
}
